# Updates the Price (D) and Volume(1h) (E) columns of the cryptos
# worksheet to the latest scraped snapshot.
#
# The Price column holds values that look numeric ("311.55") as well as
# ones that use dots as thousands separators ("28.153.75"). Excel's COM
# Value setter auto-coerces plain numeric-looking strings into real
# numbers, which would silently lose formatting/precision (e.g. "41.00"
# -> 41, "0.09600" -> 9.6E-2). Prefixing with a leading apostrophe forces
# the input to be stored as text (as it was originally), and the follow-up
# Style reset clears the "quote prefix" look so the cell keeps its
# original (unstyled) appearance.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.153.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "'1.869.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.92%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'311.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").Value = "'0.5009"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.30%  "
$ws.Range("D8").Value = "'0.3901"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.72%  "
$ws.Range("D9").Value = "'0.09600"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +24.51%  "
$ws.Range("D10").Value = "'1.138"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.05%  "
$ws.Range("D11").Value = "'41.00"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.69%  "
$ws.Range("D12").Value = "'6.459"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.19%  "
$ws.Range("D13").Value = "'20.91"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.10%  "
$ws.Range("D14").Value = "'1.874.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.10%  "
$ws.Range("D15").Value = "'1.002"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").Value = "'7.380"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.73%  "
$ws.Range("D17").Value = "'0.00001122"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.02%  "
$ws.Range("D18").Value = "'93.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.26%  "
$ws.Range("D19").Value = "'0.06617"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.73%  "
$ws.Range("D20").Value = "'17.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").Value = "'6.155"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.18%  "
$ws.Range("D23").Value = "'28.216.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.57%  "
$ws.Range("D24").Value = "'11.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.75%  "
$ws.Range("D25").Value = "'2.278"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.30%  "
$ws.Range("D26").Value = "'2.549"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.67%  "
$ws.Range("D27").Value = "'2.083.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.60%  "
$ws.Range("D28").Value = "'21.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.66%  "
$ws.Range("D29").Value = "'157.45"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.38%  "
$ws.Range("D30").Value = "'127.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("D31").Value = "'0.1054"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.32%  "
$ws.Range("D32").Value = "'1.061"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.62%  "
$ws.Range("D33").Value = "'5.623"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.63%  "
$ws.Range("D34").Value = "'3.627"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.50%  "
$ws.Range("D35").Value = "'0.06753"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.95%  "
$ws.Range("D36").Value = "'9.533"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.64%  "
$ws.Range("D37").Value = "'0.02396"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.94%  "
$ws.Range("D38").Value = "'0.2174"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.58%  "
$ws.Range("D39").Value = "'11.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.31%  "
$ws.Range("D40").Value = "'4.963"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.62%  "
$ws.Range("D41").Value = "'0.6296"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.17%  "
$ws.Range("D42").Value = "'1.176"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.41%  "
$ws.Range("D43").Value = "'1.002"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").Value = "'13.54"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.33%  "
$ws.Range("D45").Value = "'0.6022"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.42%  "
$ws.Range("D46").Value = "'3.657"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.18%  "
$ws.Range("D47").Value = "'1.258"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.55%  "
$ws.Range("D48").Value = "'123.72"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.47%  "
$ws.Range("D49").Value = "'1.978"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.82%  "
$ws.Range("D50").Value = "'1.195"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.82%  "
$ws.Range("D51").Value = "'0.06840"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.84%  "
